$d = $word.ActiveDocument

# 1) Update the arraignment date: March 05, 2022 -> March 09, 2022
$d.Content.Find.Execute(" on March 05, 2022.", $true, $false, $false, $false, $false,
                         $true, 1, $false, " on March 09, 2022.", 2)

# 2) Collapse the "by June 14, 2022 ... at 7:00 p.m" run sequence down to
#    a single run reading "by June 21, 2022".
$d.Content.Find.Execute("by June 14, 2022, and shall report to jail on June 17, 2022, at 7:00 p.m", $true, $false, $false, $false, $false,
                         $true, 1, $false, "by June 21, 2022", 2)
